# Generate Report for Handback
# Updates row 8 ("f71abb2d-5b54-4012-9ba9-746b629795cc") on both the
# zh-cn and de-de sheets with a new handback result: a Latest Target File
# hyperlink, an updated Latest Handback File name, a new Latest Handback
# DateTime, and an Error Detail message saying the handback file version
# is not the latest one. Also widens the Error Detail (P) column.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1bf45a1777fedb4e90a3bf032410776360bea26/e2e/f71abb2d-5b54-4012-9ba9-746b629795cc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8bd4e9da01da0871ec7536343d33777255c67b87/e2e/f71abb2d-5b54-4012-9ba9-746b629795cc.md."

$targetHyperlink = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b1bf45a1777fedb4e90a3bf032410776360bea26/e2e/f71abb2d-5b54-4012-9ba9-746b629795cc.md"
$targetDisplay = "f71abb2d-5b54-4012-9ba9-746b629795cc.md"

function Update-HandbackRow($SheetName, $HandbackFile, $HandbackDateTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    # I8 - Latest Target File: becomes a hyperlink to the handback markdown file
    $ws.Hyperlinks.Add($ws.Range("I8"), $targetHyperlink, "", "", $targetDisplay)
    $ws.Range("I8").Font.Underline = 2
    $ws.Range("I8").Font.Color = 15570276

    # J8 - Latest Handback File
    $ws.Range("J8").Value = $HandbackFile

    # K8 - Latest Handback DateTime
    $ws.Range("K8").Value = $HandbackDateTime

    # P8 - Error Detail
    $ws.Range("P8").Value = $errorDetail

    # Widen the Error Detail column (P, column 16) to fit the new message
    $ws.Columns.Item(16).ColumnWidth = 39.17
}

Update-HandbackRow "zh-cn" "f71abb2d-5b54-4012-9ba9-746b629795cc.9f57788fdc9b5ec9e00ee44e049d63fc4fbf56d1.zh-cn.xlf" "2016-09-05 11:02:39"
Update-HandbackRow "de-de" "f71abb2d-5b54-4012-9ba9-746b629795cc.9f57788fdc9b5ec9e00ee44e049d63fc4fbf56d1.de-de.xlf" "2016-09-05 11:02:47"

Write-Host "Handback report updated"
